$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.969.61"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "2.344.31"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'544.45"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").Value = "'134.56"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D9").Value = "2.347.02"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'5.38"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  +5.86%  "
$ws.Range("D14").Value = "2.767.95"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'23.53"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "57.998.67"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.351.44"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'10.58"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'333.40"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "'4.21"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'61.79"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'0.168"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'8.42"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "'1.40"
$ws.Range("E28").Value = "  +7.08%  "
$ws.Range("E29").Value = "  +5.00%  "
$ws.Range("D30").Value = "'170.26"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "0.0₃0728"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  +16.96%  "
$ws.Range("D34").Value = "'18.42"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'4.18"
$ws.Range("E37").Value = "  +6.27%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").Value = "'1.63"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("D40").Value = "'39.15"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "'147.23"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "'0.378"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").Value = "'286.02"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").Value = "'3.59"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'19.15"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "'0.0925"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'0.0217"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "'17.54"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").Value = "'0.380"
$ws.Range("E51").Value = "  +8.84%  "
